$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.442.78"
$ws.Range("E2").Value = "  -4.76%  "

$ws.Range("D3").Value = "3.306.72"
$ws.Range("E3").Value = "  -6.59%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'182.61"
$ws.Range("E5").Value = "  -10.61%  "

$ws.Range("D6").Value = "'532.28"
$ws.Range("E6").Value = "  -3.67%  "

$ws.Range("D7").Value = "'0.606"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "3.291.82"
$ws.Range("E8").Value = "  -6.91%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  -5.42%  "

$ws.Range("D11").Value = "'59.38"
$ws.Range("E11").Value = "  -6.36%  "

$ws.Range("E12").Value = "  -6.57%  "

$ws.Range("D13").Value = "'0.0000266"
$ws.Range("E13").Value = "  -3.12%  "

$ws.Range("D14").Value = "'9.18"
$ws.Range("E14").Value = "  -7.38%  "

$ws.Range("D15").Value = "3.841.17"
$ws.Range("E15").Value = "  -6.08%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.309.44"
$ws.Range("E16").Value = "  -6.25%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.117"
$ws.Range("E17").Value = "  -5.74%  "

$ws.Range("D18").Value = "'17.74"
$ws.Range("E18").Value = "  -5.08%  "

$ws.Range("D19").Value = "64.346.15"
$ws.Range("E19").Value = "  -4.43%  "

$ws.Range("D20").Value = "'11.15"
$ws.Range("E20").Value = "  -6.67%  "

$ws.Range("D21").Value = "'0.966"
$ws.Range("E21").Value = "  -7.06%  "

$ws.Range("D22").Value = "'375.53"
$ws.Range("E22").Value = "  -4.39%  "

$ws.Range("D23").Value = "'3.84"
$ws.Range("E23").Value = "  -5.20%  "

$ws.Range("D24").Value = "'11.27"
$ws.Range("E24").Value = "  -6.59%  "

$ws.Range("D25").Value = "'81.20"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("D26").Value = "'3.98"
$ws.Range("E26").Value = "  +5.73%  "

$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").Value = "'2.69"
$ws.Range("E28").Value = "  -4.95%  "

$ws.Range("E29").Value = "  -4.62%  "

$ws.Range("E30").Value = "  -5.36%  "

$ws.Range("D31").Value = "'29.03"
$ws.Range("E31").Value = "  -6.36%  "

$ws.Range("D32").Value = "'6.87"
$ws.Range("E32").Value = "  -6.33%  "

$ws.Range("D33").Value = "'643.12"
$ws.Range("E33").Value = "  -7.36%  "

$ws.Range("D34").Value = "'11.38"
$ws.Range("E34").Value = "  -4.00%  "

$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  -5.07%  "

$ws.Range("D36").Value = "'59.29"
$ws.Range("E36").Value = "  -7.98%  "

$ws.Range("D37").Value = "'0.399"
$ws.Range("E37").Value = "  -3.49%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'37.10"
$ws.Range("E39").Value = "  -6.33%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0730"
$ws.Range("E40").Value = "  +3.51%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.129"
$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.919.62"
$ws.Range("E43").Value = "  -5.60%  "

$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -3.17%  "

$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  -11.06%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0407"
$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +10.52%  "

$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  -2.41%  "

$ws.Range("D49").Value = "'2.62"
$ws.Range("E49").Value = "  -8.86%  "

$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("E51").Value = "  +1.04%  "
